# Generate Report for Handoff
#
# The localization status report is refreshed for a new handoff: the
# "Handed back: in sync with en-US" status becomes "Ready for handoff"
# everywhere it's used (Overview E2/F2, zh-cn C2, de-de C2), the handoff
# timestamps tied to that status move forward a minute or so, and the
# now-narrower Status column is resized to fit the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed handoff timestamps (stored as text, not numeric dates) ---
# Re-apply the datetime NumberFormat after the write so the cell keeps its
# original "yyyy-mm-dd HH:mm:ss" display style.
$overview.Range("G2").Value = "2016-09-07 05:16:46"
$overview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zhcn.Range("H2").Value = "2016-09-07 05:16:40"
$zhcn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$dede.Range("H2").Value = "2016-09-07 05:16:46"
$dede.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# --- Narrow the Status columns to fit the shorter text ---
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
